$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D51").NumberFormat = "@"

$ws.Range("D2").Value = "68.841.87"
$ws.Range("E2").Value = "  -0.58%  "
$ws.Range("D3").Value = "3.853.12"
$ws.Range("E3").Value = "  +2.50%  "
$ws.Range("E4").Value = "  -0.12%  "
$ws.Range("D5").Value = "601.21"
$ws.Range("E5").Value = "  -0.24%  "
$ws.Range("D6").Value = "161.99"
$ws.Range("E6").Value = "  -3.17%  "
$ws.Range("D7").Value = "3.850.88"
$ws.Range("E7").Value = "  +2.37%  "
$ws.Range("E8").Value = "  +0.08%  "
$ws.Range("E9").Value = "  -1.51%  "
$ws.Range("E10").Value = "  -1.24%  "
$ws.Range("E11").Value = "  -1.44%  "
$ws.Range("E12").Value = "  -0.08%  "
$ws.Range("D13").Value = "36.90"
$ws.Range("E13").Value = "  -3.04%  "
$ws.Range("E14").Value = "  -2.12%  "
$ws.Range("D15").Value = "4.491.13"
$ws.Range("E15").Value = "  +2.34%  "
$ws.Range("D16").Value = "3.836.07"
$ws.Range("E16").Value = "  +1.82%  "
$ws.Range("D17").Value = "68.984.40"
$ws.Range("E18").Value = "  +1.96%  "
$ws.Range("B19").Value = "Uniswap"
$ws.Range("C19").Value = "https://coinranking.com/coin/_H5FVG9iW+uniswap-uni"
$ws.Range("D19").Value = "11.49"
$ws.Range("E19").Value = "  +1.65%  "
$ws.Range("B20").Value = "TRON"
$ws.Range("C20").Value = "https://coinranking.com/coin/qUhEFk1I61atv+tron-trx"
$ws.Range("D20").Value = "0.113"
$ws.Range("E20").Value = "  -0.20%  "
$ws.Range("D21").Value = "17.13"
$ws.Range("E21").Value = "  -1.11%  "
$ws.Range("D22").Value = "483.63"
$ws.Range("E22").Value = "  -2.06%  "
$ws.Range("D23").Value = "0.719"
$ws.Range("E23").Value = "  -1.39%  "
$ws.Range("D24").Value = "0.0000158"
$ws.Range("E24").Value = "  +3.75%  "
$ws.Range("D25").Value = "83.91"
$ws.Range("E25").Value = "  -1.15%  "
$ws.Range("E26").Value = "  -2.75%  "
$ws.Range("E27").Value = "  -1.74%  "
$ws.Range("D28").Value = "0.999"
$ws.Range("E28").Value = "  -0.06%  "
$ws.Range("D29").Value = "9.99"
$ws.Range("E29").Value = "  -1.29%  "
$ws.Range("E30").Value = "  -0.97%  "
$ws.Range("D31").Value = "7.93"
$ws.Range("E31").Value = "  -2.61%  "
$ws.Range("D32").Value = "3.995.76"
$ws.Range("E32").Value = "  +2.35%  "
$ws.Range("E33").Value = "  -4.34%  "
$ws.Range("D34").Value = "32.20"
$ws.Range("E34").Value = "  +1.86%  "
$ws.Range("D35").Value = "3.794.89"
$ws.Range("E35").Value = "  +2.65%  "
$ws.Range("E36").Value = "  -1.42%  "
$ws.Range("E37").Value = "  +1.17%  "
$ws.Range("D38").Value = "0.140"
$ws.Range("E38").Value = "  +2.41%  "
$ws.Range("D39").Value = "5.89"
$ws.Range("E39").Value = "  -2.12%  "
$ws.Range("E40").Value = "  -0.09%  "
$ws.Range("E41").Value = "  -1.94%  "
$ws.Range("D42").Value = "438.20"
$ws.Range("E42").Value = "  +1.68%  "
$ws.Range("D43").Value = "2.97"
$ws.Range("E43").Value = "  -1.48%  "
$ws.Range("E44").Value = "  -0.47%  "
$ws.Range("D45").Value = "1.98"
$ws.Range("E45").Value = "  -0.58%  "
$ws.Range("D47").Value = "8.39"
$ws.Range("E47").Value = "  -0.98%  "
$ws.Range("D48").Value = "143.29"
$ws.Range("E48").Value = "  +1.86%  "
$ws.Range("D49").Value = "2.834.22"
$ws.Range("E49").Value = "  +1.30%  "
$ws.Range("E50").Value = "  +1.59%  "
$ws.Range("D51").Value = "25.71"
$ws.Range("E51").Value = "  +11.31%  "
